# Iraq League base update (14-04-2024 18:28)
#
# 1) Eight pairs of adjacent match rows had their match-specific data
#    (id/B, HomeTeam/F, AwayTeam/G, FTHG/H, FTAG/I, FTR/J, and all the
#    odds columns K:AC) swapped between the two rows in the pair. The
#    row-sequence column A and the Div/Div Original Name/Date columns
#    (C, D, E) are identical within each pair, so they are left as-is.
# 2) One new match (row 160) was appended at the end of the sheet - a
#    fixture that hasn't kicked off yet, so it has no FTHG/FTAG/FTR or
#    PL_AhOver/PL_AhUnder values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    # Columns B, F:AC (skip A, C, D, E) -> col indices 2 and 6..29
    $colIdxs = @(2) + (6..29)
    foreach ($col in $colIdxs) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

$pairs = @(
    @(17, 18),
    @(22, 23),
    @(78, 79),
    @(91, 92),
    @(103, 104),
    @(108, 109),
    @(135, 136),
    @(151, 152)
)

foreach ($pair in $pairs) {
    Swap-RowData $pair[0] $pair[1]
}

# New row 160 - upcoming fixture, not yet played.
$r = 160

$a = $ws.Cells.Item($r, 1)
$a.Value2 = 158
$a.Font.Bold = $true
$a.Borders.LineStyle = 1
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160

$ws.Cells.Item($r, 2).Value2 = 8088065
$ws.Cells.Item($r, 3).Value2 = "Iraq League"
$ws.Cells.Item($r, 4).Value2 = "Iraq League"

$e = $ws.Cells.Item($r, 5)
$e.Value2 = 45397.66666666666
$e.NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($r, 6).Value2 = "Al Quwa Al Jawiya"
$ws.Cells.Item($r, 7).Value2 = "Naft Maysan"

# H, I, J (FTHG/FTAG/FTR) intentionally left blank - match not played yet.

$ws.Cells.Item($r, 11).Value2 = 1.6
$ws.Cells.Item($r, 12).Value2 = 3.5
$ws.Cells.Item($r, 13).Value2 = 5
$ws.Cells.Item($r, 14).Value2 = 1.444
$ws.Cells.Item($r, 15).Value2 = 4
$ws.Cells.Item($r, 16).Value2 = 6
$ws.Cells.Item($r, 17).Value2 = -1.25
$ws.Cells.Item($r, 18).Value2 = 1.95
$ws.Cells.Item($r, 19).Value2 = 1.85
$ws.Cells.Item($r, 20).Value2 = 2.75
$ws.Cells.Item($r, 21).Value2 = 1.825
$ws.Cells.Item($r, 22).Value2 = 1.975
$ws.Cells.Item($r, 23).Value2 = 0
$ws.Cells.Item($r, 24).Value2 = 0
$ws.Cells.Item($r, 25).Value2 = 0
$ws.Cells.Item($r, 26).Value2 = 0
$ws.Cells.Item($r, 27).Value2 = 0

# AB, AC (PL_AhOver/PL_AhUnder) intentionally left blank - match not played yet.
